$d = $word.ActiveDocument

function Set-ParagraphXml($para, [string]$xmlPackage) {
    $r = $para.Range
    $r.End = $r.End - 1
    $r.Text = ""
    $r2 = $para.Range
    $r2.End = $r2.End - 1
    $r2.InsertXML($xmlPackage)
}

# Locate the three target paragraphs by their distinctive leading text.
$progParaPt = $null
$progParaEn = $null
$bibPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($progParaPt -eq $null -and $t.StartsWith("- Breve hist")) {
        $progParaPt = $p
    } elseif ($progParaEn -eq $null -and $t.StartsWith("- Brief history")) {
        $progParaEn = $p
    } elseif ($bibPara -eq $null -and $t.StartsWith("Artaxo P (2014)")) {
        $bibPara = $p
    }
}

if ($progParaPt -eq $null) { throw "Portuguese Programa paragraph not found" }
if ($progParaEn -eq $null) { throw "English Programa paragraph not found" }
if ($bibPara -eq $null) { throw "Bibliografia paragraph not found" }

$ptXml = '<?xml version="1.0" encoding="UTF-16" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>- Breve história da Terra, com ênfase na formação e evolução da atmosfera terrestre;</w:t><w:br/><w:t>- Mudanças climáticas na história da Terra (escalas de tempo geológica e ecológica)</w:t><w:br/><w:t>- Sistemas atmosfera-hidrosfera-criosfera-biosfera-litosfera</w:t><w:br/><w:t xml:space="preserve">- Ciclos biogeoquímicos e marcadores isotópicos de mudanças climáticas (ciclo da água e do carbono) </w:t><w:br/><w:t xml:space="preserve">- Conceitos e histórico das mudanças climáticas antropogênicas (impactos humanos, revolução industrial e atualidade), </w:t><w:br/><w:t>- O Antropoceno e os principais drivers das mudanças climáticas antropogênicas no mundo (setores de energia, construção, indústria, transporte) e no Brasil (agricultura, desmatamento e urbanização);</w:t><w:br/><w:t>- Eventos climáticos extremos, ordenamento territorial, riscos geológicos e saúde planetária (ênfase em áreas altamente urbanizadas);</w:t><w:br/><w:t xml:space="preserve">- Estratégias de adaptação e mitigação das mudanças climáticas (Soluções baseadas na Natureza – SbN, Florestas Urbanas); </w:t><w:br/><w:t>- O debate nacional e internacional sobre as mudanças climáticas antropogênicas. Acordos internacionais. O surgimento do Painel Intergovernamental sobre Mudanças Climáticas (IPCC), relatórios e evolução dos modelos e cenários previstos. A Agenda Climática e de Sustentabilidade (Agenda 2030);</w:t><w:br/><w:t>- Políticas públicas e diretrizes internacionais sobre mudanças climáticas. Contribuições Nacionalmente Determinadas (iNDC). Propostas e acordos locais, nacionais e internacionais para alcançar a sustentabilidade global.</w:t><w:br/><w:t xml:space="preserve">- Inventário de emissões de gases de efeito estufa (Protocolo GHG/FGV) – atividade prática dos escopos 1 e 2; </w:t><w:br/><w:t>- Saída de campo (Centro Nacional de Monitoramento e Alertas de Desastres Naturais - CEMADEM). Atividades direcionadas e tratamento de dados atualizados.</w:t><w:br/><w:t>- Elaboração e apresentação de projeto sobre adaptação de cidades à futuros extremos climáticos (selecionar uma cidade do Vale do Paraíba, analisar sistematicamente o território no âmbito das vulnerabilidades climática e socioambiental e propor mecanismos de adaptação da cidade para tais cenários previstos do clima para as próximas décadas).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$enXml = '<?xml version="1.0" encoding="UTF-16" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:i/></w:rPr><w:t>- Brief history of the Earth, with emphasis on the formation and evolution of the Earth''s atmosphere;</w:t><w:br/><w:t>- Climate change in Earth''s history (geological and ecological timescales)</w:t><w:br/><w:t>- Atmosphere-hydrosphere-cryosphere-biosphere-lithosphere systems</w:t><w:br/><w:t>- Biogeochemical cycles and isotopic markers of climate change (water and carbon cycle)</w:t><w:br/><w:t>- Concepts and history of anthropogenic climate change (human impacts, industrial revolution and current affairs),</w:t><w:br/><w:t>- The Anthropocene and the main drivers of anthropogenic climate change in the world (energy, construction, industry, transport sectors) and in Brazil (agriculture, deforestation and urbanization);</w:t><w:br/><w:t>- Extreme climate events, territorial planning, geological risks and planetary health (emphasis on highly urbanized areas);</w:t><w:br/><w:t>- Climate change adaptation and mitigation strategies (Nature-based Solutions – SbN, Urban Forests);</w:t><w:br/><w:t>- The national and international debate on anthropogenic climate change. International agreements. The emergence of the Intergovernmental Panel on Climate Change (IPCC), reports and evolution of predicted models and scenarios. The Climate and Sustainability Agenda (Agenda 2030);</w:t><w:br/><w:t>- Public policies and international guidelines on climate change. Nationally Determined Contributions (iNDC). Local, national and international proposals and agreements to achieve global sustainability.</w:t><w:br/><w:t>- Inventory of greenhouse gas emissions (GHG/FGV Protocol) – practical activity of scopes 1 and 2;</w:t><w:br/><w:t>- Field trip (National Center for Natural Disaster Monitoring and Alerts - CEMADEM). Targeted activities and updated data processing.</w:t><w:br/><w:t>- Preparation and presentation of a project on adapting cities to future climate extremes (select a city in the Paraíba Valley, systematically analyze the territory in terms of climate and socio-environmental vulnerabilities and propose adaptation mechanisms for the city to such predicted climate scenarios for the next decades).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$bibXml = '<?xml version="1.0" encoding="UTF-16" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Artaxo P (2014). Uma nova era geológica em nosso planeta: o Antropoceno?. Revista Usp, (103), 13-24. </w:t><w:br/><w:br/><w:t>Iwama AY, Batistella M, &amp; Ferreira LDC (2014). Riscos geotécnicos e vulnerabilidade social em zonas costeiras: desigualdades e mudanças climáticas. Ambiente &amp; Sociedade, 17, 251-274.</w:t><w:br/><w:br/><w:t>Macreadie PI, Costa MD, Atwood TB, Friess DA, Kelleway JJ, Kennedy H, ... &amp; Duarte CM (2021). Blue carbon as a natural climate solution. Nature Reviews Earth &amp; Environment, 2(12), 826-839.</w:t><w:br/><w:br/><w:t>Mitchard ET (2018). The tropical forest carbon cycle and climate change. Nature, 559(7715), 527-534.</w:t><w:br/><w:br/><w:t xml:space="preserve">Novello VF, da Cruz FW, Vuille M, Campos JLPS, Stríkis NM, Apaéstegui J, ... &amp; Karmann I (2021). Investigating δ13C values in stalagmites from tropical South America for the last two millennia. Quaternary Science Reviews, 255, 106822. </w:t><w:br/><w:br/><w:t xml:space="preserve">Oliveira MJ, Carneiro CDR, da Silva Vecchia FA, &amp; de Mello Baptista GM (2017). Ciclos climáticos e causas naturais das mudanças do clima. Terrae didática, 13(3), 149-184. </w:t><w:br/><w:br/><w:t>Oki T, Entekhabi D, &amp; Harrold TI (1999). The global water cycle. Global energy and water cycles, 10, 27.</w:t><w:br/><w:br/><w:t>Pereira P, Wang F, Inacio M, Kalinauskas, M, Bogdzevič K, Bogunovic I, ... &amp; Barcelo D (2024). Nature-based solutions for carbon sequestration in urban environments. Current Opinion in Environmental Science &amp; Health, 100536.</w:t><w:br/><w:br/><w:t>Romanello M, Di Napoli C, Green C, Kennard H, Lampard P, Scamman D, ... &amp; Costello A (2023). The 2023 report of the Lancet Countdown on health and climate change: the imperative for a health-centred response in a world facing irreversible harms. The Lancet, 402(10419), 2346-2394.</w:t><w:br/><w:br/><w:t>Seddon N, Smith A, Smith P, Key I, Chausson A, Girardin C, ... &amp; Turner B (2021). Getting the message right on nature‐based solutions to climate change. Global change biology, 27(8), 1518-1546.</w:t><w:br/><w:br/><w:t>Waters CN, Zalasiewicz J, Summerhayes C, Barnosky AD, Poirier C, Gałuszka A, ... &amp; Wolfe AP (2016). The Anthropocene is functionally and stratigraphically distinct from the Holocene. Science, 351(6269), aad2622.</w:t><w:br/><w:br/><w:t>Diversos artigos científicos na Revista “Global Change Biology”, Grupo Wiley.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

Set-ParagraphXml $progParaPt $ptXml
Set-ParagraphXml $progParaEn $enXml
Set-ParagraphXml $bibPara $bibXml

Write-Host "Done editing Programa (PT/EN) and Bibliografia paragraphs."
